$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1750
$ws.Range("I111").Value = 1500
$ws.Range("K111").Value = 4500
$ws.Range("M111").Value = -1433
$ws.Range("H113").Value = 3517.5
$ws.Range("I113").Value = 3444.8
$ws.Range("K113").Value = 3444.8
$ws.Range("M113").Value = -190.8000000000002
$ws.Range("H137").Value = 3580034.8
$ws.Range("I137").Value = 5760.905
$ws.Range("K137").Value = 17282.715
$ws.Range("M137").Value = -14732.715
$ws.Range("H138").Value = 5583.5684
$ws.Range("I138").Value = 12349.857
$ws.Range("J138").Value = 2425.9666
$ws.Range("K138").Value = 37049.571
$ws.Range("L138").Value = 7277.899800000001
$ws.Range("M138").Value = -31909.571
$ws.Range("N138").Value = -17557.8998

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 92499.75
$ws.Range("J24").Value = 92499.75
$ws.Range("L24").Value = 92499.75
$ws.Range("N24").Value = -93247.75
$ws.Range("H32").Value = 297493.06
$ws.Range("I32").Value = 385221.34
$ws.Range("K32").Value = 385221.34
$ws.Range("M32").Value = -384934.34
$ws.Range("H45").Value = 47590.047
$ws.Range("I45").Value = 54708.79
$ws.Range("J45").Value = 2504.6667
$ws.Range("K45").Value = 54708.79
$ws.Range("L45").Value = 2504.6667
$ws.Range("M45").Value = -54331.79
$ws.Range("N45").Value = -3258.6667
$ws.Range("H80").Value = 49099
$ws.Range("J80").Value = 49099
$ws.Range("L80").Value = 49099
$ws.Range("N80").Value = -51095
$ws.Range("H83").Value = 49099
$ws.Range("J83").Value = 49099
$ws.Range("L83").Value = 147297
$ws.Range("N83").Value = -157281
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H100").Value = 92499.75
$ws.Range("J100").Value = 92499.75
$ws.Range("L100").Value = 92499.75
$ws.Range("N100").Value = -94663.75
$ws.Range("H102").Value = 2855
$ws.Range("I102").Value = 2855
$ws.Range("K102").Value = 2855
$ws.Range("M102").Value = -1233
$ws.Range("H122").Value = 3333.3333
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 4250
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 12750
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -17650
$ws.Range("H132").Value = 3526.3572
$ws.Range("I132").Value = 2516.0908
$ws.Range("J132").Value = 4180.0586
$ws.Range("K132").Value = 7548.2724
$ws.Range("L132").Value = 12540.1758
$ws.Range("M132").Value = -5018.2724
$ws.Range("N132").Value = -17600.1758

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 80000
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 80000
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -81148
$ws.Range("H105").Value = 41252.25
$ws.Range("I105").Value = 36669.668
$ws.Range("K105").Value = 36669.668
$ws.Range("M105").Value = -34922.668
$ws.Range("H107").Value = 34329.5
$ws.Range("I107").Value = 40595.4
$ws.Range("K107").Value = 40595.4
$ws.Range("M107").Value = -38675.4
$ws.Range("H134").Value = 28126770
$ws.Range("I134").Value = 1700.909
$ws.Range("J134").Value = 90001930
$ws.Range("K134").Value = 5102.727000000001
$ws.Range("L134").Value = 270005790
$ws.Range("M134").Value = -2567.727000000001
$ws.Range("N134").Value = -270010860

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 221.16667
$ws.Range("J7").Value = 56.25
$ws.Range("L7").Value = 56.25
$ws.Range("N7").Value = -282.25
$ws.Range("H22").Value = 1191.5
$ws.Range("J22").Value = 588.3333
$ws.Range("L22").Value = 588.3333
$ws.Range("N22").Value = -1288.3333
$ws.Range("H99").Value = 10001050
$ws.Range("I99").Value = 10001050
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 10001050
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 3010
$ws.Range("I105").Value = 2015.125
$ws.Range("K105").Value = 2015.125
$ws.Range("M105").Value = -268.125
$ws.Range("H126").Value = 10001050
$ws.Range("I126").Value = 10001050
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 30003150
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 2942.3333
$ws.Range("I132").Value = 3975
$ws.Range("J132").Value = 2426
$ws.Range("K132").Value = 11925
$ws.Range("L132").Value = 7278
$ws.Range("M132").Value = -9395
$ws.Range("N132").Value = -12338
$ws.Range("H134").Value = 3326.7693
$ws.Range("I134").Value = 4311.5
$ws.Range("K134").Value = 12934.5
$ws.Range("M134").Value = -10399.5
$ws.Range("H141").Value = 264166.34
$ws.Range("J141").Value = 264166.34
$ws.Range("L141").Value = 264166.34
$ws.Range("N141").Value = -274526.34

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1154.4546
$ws.Range("J17").Value = 1169.9
$ws.Range("L17").Value = 3509.7
$ws.Range("N17").Value = -3847.7
$ws.Range("H34").Value = 1300
$ws.Range("J34").Value = 1750
$ws.Range("L34").Value = 5250
$ws.Range("N34").Value = -5418
$ws.Range("H39").Value = 7207.3335
$ws.Range("I39").Value = 899.3333
$ws.Range("J39").Value = 10361.333
$ws.Range("K39").Value = 2697.9999
$ws.Range("L39").Value = 31083.999
$ws.Range("M39").Value = -2403.9999
$ws.Range("N39").Value = -31671.999
$ws.Range("H50").Value = 3675.75
$ws.Range("J50").Value = 3801.1428
$ws.Range("L50").Value = 11403.4284
$ws.Range("N50").Value = -12365.4284
$ws.Range("H53").Value = 3675.75
$ws.Range("J53").Value = 3801.1428
$ws.Range("L53").Value = 11403.4284
$ws.Range("N53").Value = -12365.4284
$ws.Range("H55").Value = 5226
$ws.Range("J55").Value = 6248.25
$ws.Range("L55").Value = 18744.75
$ws.Range("N55").Value = -19098.75
$ws.Range("H57").Value = 8476.25
$ws.Range("I57").Value = 4452.5
$ws.Range("K57").Value = 13357.5
$ws.Range("M57").Value = -12798.5
$ws.Range("H99").Value = 111123016
$ws.Range("J99").Value = 15867.333
$ws.Range("L99").Value = 47601.999
$ws.Range("N99").Value = -52093.999
$ws.Range("H117").Value = 60608836
$ws.Range("J117").Value = 60608836
$ws.Range("L117").Value = 181826508
$ws.Range("N117").Value = -181833392

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 97999.5
$ws.Range("J86").Value = 97999.5
$ws.Range("L86").Value = 97999.5
$ws.Range("N86").Value = -100371.5
$ws.Range("H89").Value = 97999.5
$ws.Range("J89").Value = 97999.5
$ws.Range("L89").Value = 293998.5
$ws.Range("N89").Value = -305854.5
$ws.Range("H95").Value = 47649
$ws.Range("J95").Value = 47649
$ws.Range("L95").Value = 47649
$ws.Range("N95").Value = -53141
$ws.Range("H102").Value = 16130349
$ws.Range("I102").Value = 20001180
$ws.Range("J102").Value = 1889
$ws.Range("K102").Value = 20001180
$ws.Range("L102").Value = 1889
$ws.Range("M102").Value = -19999558
$ws.Range("N102").Value = -5133
$ws.Range("H122").Value = 3234.3809
$ws.Range("I122").Value = 3055
$ws.Range("K122").Value = 9165
$ws.Range("M122").Value = -6715
$ws.Range("H132").Value = 1071130
$ws.Range("I132").Value = 1100.3334
$ws.Range("K132").Value = 3301.0002
$ws.Range("M132").Value = -771.0001999999999
$ws.Range("H135").Value = 75000
$ws.Range("J135").Value = 75000
$ws.Range("L135").Value = 75000
$ws.Range("N135").Value = -85140

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3069.611
$ws.Range("I122").Value = 2904.8965
$ws.Range("J122").Value = 3752
$ws.Range("K122").Value = 8714.6895
$ws.Range("L122").Value = 11256
$ws.Range("M122").Value = -6264.6895
$ws.Range("N122").Value = -16156
$ws.Range("H132").Value = 4384.7856
$ws.Range("I132").Value = 2399.6667
$ws.Range("K132").Value = 7199.000100000001
$ws.Range("M132").Value = -4669.000100000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1787525.9
$ws.Range("I107").Value = 1423.375
$ws.Range("J107").Value = 3573628.5
$ws.Range("K107").Value = 4270.125
$ws.Range("L107").Value = 10720885.5
$ws.Range("M107").Value = -2350.125
$ws.Range("N107").Value = -10724725.5
$ws.Range("H132").Value = 1950
$ws.Range("I132").Value = 1631.9131
$ws.Range("K132").Value = 4895.7393
$ws.Range("M132").Value = -2365.7393
